$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.561.83"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "'1.573.03"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  -1.23%  "
$ws.Range("D5").Value = "'211.68"
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("D6").Value = "'0.493"
$ws.Range("E6").Value = "  +0.76%  "
$ws.Range("D7").Value = "'0.991"
$ws.Range("E7").Value = "  -1.43%  "
$ws.Range("D8").Value = "'22.57"
$ws.Range("E8").Value = "  +2.27%  "
$ws.Range("D9").Value = "'0.252"
$ws.Range("E9").Value = "  +1.09%  "
$ws.Range("D10").Value = "'0.0598"
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("D12").Value = "'1.798.16"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").Value = "'1.564.28"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").Value = "'27.534.15"
$ws.Range("E16").Value = "  +1.85%  "
$ws.Range("D17").Value = "'62.09"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "'225.82"
$ws.Range("E18").Value = "  +4.74%  "
$ws.Range("D19").Value = "'7.57"
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "'0.992"
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("D22").Value = "'4.18"
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("D23").Value = "'9.45"
$ws.Range("E23").Value = "  +2.61%  "
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("D25").Value = "'150.77"
$ws.Range("E25").Value = "  -1.52%  "
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").Value = "'15.23"
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("E28").Value = "  +2.17%  "
$ws.Range("D29").Value = "'0.993"
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("D30").Value = "'1.14"
$ws.Range("E30").Value = "  +1.64%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "'3.26"
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("D34").Value = "'1.457.25"
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("E35").Value = "  +4.67%  "
$ws.Range("E36").Value = "  +1.98%  "
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("E39").Value = "  +1.58%  "
$ws.Range("D40").Value = "'0.820"
$ws.Range("E40").Value = "  +1.32%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'5.82"
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.35"
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("D43").Value = "'0.993"
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("D44").Value = "'65.60"
$ws.Range("E44").Value = "  +1.27%  "
$ws.Range("D45").Value = "'0.969"
$ws.Range("E45").Value = "  -3.12%  "
$ws.Range("E46").Value = "  +2.55%  "
$ws.Range("D47").Value = "'1.711.54"
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("D48").Value = "'86.61"
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("D50").Value = "'0.0₇0965"
$ws.Range("E50").Value = "  -6.89%  "
$ws.Range("E51").Value = "  -1.88%  "
